# Applies the "Add files via upload" revision to Proposal.docx
# (replay-value / procedural-generation proposal) using Word COM-interop
# Find & Replace calls against $word.ActiveDocument.

$d = $word.ActiveDocument

# Problem Statement paragraph
$d.Content.Find.Execute("Video games lack replay value and to solve this issue procedurally generated worlds, multiple starting classes, and numerous secrets should be implemented in", $true, $false, $false, $false, $false, $true, 1, $false, "Video games lack replay value and a way to solve this issue is with procedurally generated worlds, customization, and numerous secrets should be implemented in", 2) | Out-Null

# "To address the problem..." paragraph
$d.Content.Find.Execute("procedural world generation, multiple starting classes, and secrets for the player to find. This will allow players", $true, $false, $false, $false, $false, $true, 1, $false, "procedural world generation and multiple weapon types for customization. This will allow players", 2) | Out-Null

# "In conclusion..." paragraph
$d.Content.Find.Execute("multiple starting classes for the player to choose from, and numerous secrets", $true, $false, $false, $false, $false, $true, 1, $false, "multiple weapon types for customization, and numerous secrets", 2) | Out-Null

# Project Description paragraph - bosses
$d.Content.Find.Execute("Because of procedural generating the bosses will be randomly spread out so the player can fight the bosses in whatever order they would like.", $true, $false, $false, $false, $false, $true, 1, $false, "Bosses will be randomly spawned throughout the map so the player can fight the bosses in whatever order they would like.", 2) | Out-Null

# Project Description paragraph - enemies distribution / remove starting classes
$d.Content.Find.Execute("Enemies, plants, and merchants will be randomly distributed throughout the world but will be confined to certain regions. Also, there will be multiple starting character classes for the player to choose from.", $true, $false, $false, $false, $false, $true, 1, $false, "Enemies, checkpoints, and foliage will be randomly distributed throughout the world.", 2) | Out-Null

# Project Description paragraph - loot drops
$d.Content.Find.Execute("Each class will have different stats and different perk weapons that will deal more damage. Enemies that the player defeats will drop different items of random quality based on the enemy’s tier. The player will be able to loot these items.", $true, $false, $false, $false, $false, $true, 1, $false, "Enemies that the player defeats will drop different weapons, armor, and arrows. The player will be able to loot these items.", 2) | Out-Null

# Project Description paragraph - remove crafting/merchants sentence
$d.Content.Find.Execute("Also, the player will be able to harvest materials from plants and chest that they can then use to craft different items like health potions. The player will be able to visit merchants to buy weapons, armor, and crafting materials. Furthermore, c", $true, $false, $false, $false, $false, $true, 1, $false, "Furthermore, c", 2) | Out-Null

# Project Description paragraph - checkpoints/respawn rewrite
$d.Content.Find.Execute("If the player rest at the checkpoint, then their health and health potions will be restored but resting will respawn enemies. The last checkpoint that the player visited will be set to the player’s respawn point for if an enemy kills them. Finally, the player will have a map of the land that they have explored, marking the player’s current location, merchant locations, discovered boss locations, and visited checkpoints. The player will be able to fast travel between visited checkpoints by using the map.", $true, $false, $false, $false, $false, $true, 1, $false, "If the player rests at the checkpoint, then their health will be restored. The last checkpoint that the player visited will be set as the player’s respawn point. So, if the player is killed then they will respawn at the last checkpoint they interacted with; if they did not interact with a  checkpoint then they will respawn at the starting spawn point.", 2) | Out-Null

# Additional Software/Equipment Needed paragraph (first occurrence only)
$d.Content.Find.Execute("Unreal Engine 5, Visual Studio, and Blender.", $true, $false, $false, $false, $false, $true, 1, $false, "Unreal Engine 5", 1) | Out-Null

# Outline of Future Research Efforts paragraph
$d.Content.Find.Execute("Future research efforts will include but are not limited to downloading and installing Unreal Engine 5, Visual Studio, and Blender. ", $true, $false, $false, $false, $false, $true, 1, $false, "Future research efforts will include but are not limited to downloading and installing Unreal Engine 5. ", 2) | Out-Null

Write-Output "done"
